$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new pin-name labels in column A ---
# Values are assigned in the same order they were first typed in the
# original edit so the shared-string table comes out in the matching order.
$ws.Range("A6").Value = "Eng 1 LP 2"
$ws.Range("A8").Value = "Eng 1 LP 3"
$ws.Range("A16").Value = "Eng 2 LP 1"
$ws.Range("A18").Value = "Eng 2 LP 2"
$ws.Range("A7").Value = "Eng 2 HP 1"
$ws.Range("A9").Value = "Eng 2 HP 2"
$ws.Range("A11").Value = "Eng 2 HP 3"
$ws.Range("A10").Value = "Eng 1 HP 1"
$ws.Range("A12").Value = "Eng 1 HP 2"
$ws.Range("A14").Value = "Eng 1 HP 3"
$ws.Range("A15").Value = "Eng 1 TGT -"
$ws.Range("A13").Value = "Eng 1 TGT +"
$ws.Range("A19").Value = "Eng 2 TGT -"
$ws.Range("A17").Value = "Eng 2 TGT +"
$ws.Range("A32").Value = "Eng 1 Throt +"
$ws.Range("A34").Value = "Eng 1 Throt -"
$ws.Range("A36").Value = "Eng 1 Throt out"
$ws.Range("A40").Value = "Eng 2 Throt -"
$ws.Range("A38").Value = "Eng 2 Throt +"
$ws.Range("A42").Value = "Eng 2 Throt out"
$ws.Range("A46").Value = "Air start"
$ws.Range("A33").Value = "Eng 1 Master"
$ws.Range("A37").Value = "Eng 2 Master"
$ws.Range("A41").Value = "Eng 1 Fuel cock"
$ws.Range("A45").Value = "Eng 2 fuel cock"
$ws.Range("A47").Value = "5V, air start"
$ws.Range("A44").Value = "5V, Eng 2 fuel cock"
$ws.Range("A43").Value = "5V, Eng 1 Fuel cock"
$ws.Range("A35").Value = "5V, Eng 1 Master"
$ws.Range("A39").Value = "5V, Eng 2 Master"
$ws.Range("A26").Value = "Eng 2 oil press 1"
$ws.Range("A28").Value = "Eng 2 oil press 2"
$ws.Range("A30").Value = "Eng 2 oil press 3"
$ws.Range("A20").Value = "Eng 1 oil press 1"
$ws.Range("A22").Value = "Eng 1 oil press 2"
$ws.Range("A24").Value = "Eng 1 oil press 3"
$ws.Range("A21").Value = "LP spin +"
$ws.Range("A23").Value = "LP spin - (gap)"
$ws.Range("A25").Value = "Eng 1 start"
$ws.Range("A29").Value = "Eng 2 start"
$ws.Range("A31").Value = "5V, Eng 2 start"
$ws.Range("A27").Value = "5V, Eng 1 start"

# --- Apply explicit black font colour to the "noted bug" rows ---
# (creates the 2nd font in styles.xml: Calibri 12 rgb(000000))
$ws.Range("A7:A17").Font.Color = 0
$ws.Range("A19:A31").Font.Color = 0
$ws.Range("A39").Font.Color = 0

# --- Bold the header row ---
# (creates the 3rd font in styles.xml: Calibri 12 bold, theme colour)
$ws.Range("A3:E3").Font.Bold = $true

# --- Column A width ---
$ws.Columns("A").ColumnWidth = 16.166666666666668

# --- Update selection to match the saved view ---
$ws.Range("A28").Select()
